# Fruta / hortaliza, semanal
#
# The underlying weekly data pull reshuffled which calendar week each
# "Comercializadora del Agro de Limari - Tuna" price quote belongs to.
# Every (Fecha + Calidad/Volumen/Precios) record is still intact - whole
# row-groups (one per Fecha) were simply relocated to different row
# ranges in the sheet. Rebuild the data block (A2:T38) under the new
# row order and write it back in a single pass so no row is read after
# it has already been overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A2:T38")
$src = $srcRange.Value()          # 1-based COM SafeArray: $src[row, col], row=1..37, col=1..20

$numRows = $src.GetLength(0)
$numCols = $src.GetLength(1)

# Maps the 0-based row offset of the OUTPUT block (within A2:T38) to the
# 1-based row index of the INPUT block ($src) it is copied from.
$rowMap = @{
    0 = 33
    1 = 34
    2 = 35
    3 = 17
    4 = 18
    5 = 19
    6 = 20
    7 = 21
    8 = 22
    9 = 23
    10 = 4
    11 = 5
    12 = 6
    13 = 30
    14 = 31
    15 = 32
    16 = 24
    17 = 25
    18 = 26
    19 = 27
    20 = 28
    21 = 29
    22 = 1
    23 = 2
    24 = 3
    25 = 7
    26 = 8
    27 = 9
    28 = 10
    29 = 14
    30 = 15
    31 = 16
    32 = 36
    33 = 37
    34 = 11
    35 = 12
    36 = 13
}

$dst = New-Object 'object[,]' $numRows, $numCols   # 0-based .NET array: $dst[row, col], row/col = 0..n-1

for ($r = 0; $r -lt $numRows; $r++) {
    $fromRow = $rowMap[$r]
    for ($c = 1; $c -le $numCols; $c++) {
        $dst[$r, $c - 1] = $src[$fromRow, $c]
    }
}

$srcRange.Value = $dst
